# Slide 5 ("Updates Since IETF-106 (Revision-00)"), Content Placeholder:
# insert a new bullet "Addressed review comments" right after the existing
# "Added procedure for hop-by-hop IOAM" bullet (same lvl=1 / Wingdings "ü"
# bullet), pushing "Various editorial changes" down by one paragraph.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Paragraph 2 is "Added procedure for hop-by-hop IOAM" - it already carries
# the exact paragraph formatting (lvl 1, Wingdings "ü" bullet) and run
# formatting (lang="en-CA") that the new bullet needs to match. Inserting a
# new paragraph *before* it clones that same formatting onto the new
# paragraph (inserting after/before the following paragraph would instead
# clone that other paragraph's formatting, e.g. lang="en-US").
$added = $tr.Paragraphs(2, 1)
$added.InsertBefore("Addressed review comments`r") | Out-Null

# The clone now sits where "Added procedure..." used to be (new paragraph 2)
# and the original "Added procedure..." text got pushed to paragraph 3.
# Swap the two runs' text back into the order required by the edit, which
# keeps each paragraph's own (already-correct) formatting untouched.
$tr = $sh.TextFrame.TextRange
$p2 = $tr.Paragraphs(2, 1)
$p3 = $tr.Paragraphs(3, 1)

$p2.Runs(1, 1).Text = "Added procedure for hop-by-hop IOAM"
$p3.Runs(1, 1).Text = "Addressed review comments"
